$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the data range is treated as text so values like "1.000" or "30.20"
# are not silently reinterpreted as numbers and lose formatting (e.g. trailing zeros).
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range('B2').Value = 'Bitcoin'
$ws.Range('C2').Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$ws.Range('D2').Value = '26.057.34'
$ws.Range('E2').Value = '  +1.56%  '
$ws.Range('B3').Value = 'Ethereum'
$ws.Range('C3').Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$ws.Range('D3').Value = '1.764.31'
$ws.Range('E3').Value = '  +1.32%  '
$ws.Range('B4').Value = 'TetherUSD'
$ws.Range('C4').Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$ws.Range('D4').Value = '0.9990'
$ws.Range('E4').Value = '  -0.35%  '
$ws.Range('B5').Value = 'BNB'
$ws.Range('C5').Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D5').Value = '237.85'
$ws.Range('E5').Value = '  -0.08%  '
$ws.Range('B6').Value = 'USDC'
$ws.Range('C6').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D6').Value = '0.9993'
$ws.Range('E6').Value = '  -0.25%  '
$ws.Range('B7').Value = 'XRP'
$ws.Range('C7').Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D7').Value = '0.5252'
$ws.Range('E7').Value = '  +4.47%  '
$ws.Range('B8').Value = 'Cardano'
$ws.Range('C8').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D8').Value = '0.2725'
$ws.Range('E8').Value = '  +3.19%  '
$ws.Range('B9').Value = 'OKB'
$ws.Range('C9').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D9').Value = '40.43'
$ws.Range('E9').Value = '  -3.07%  '
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').Value = '0.06221'
$ws.Range('E10').Value = '  +2.16%  '
$ws.Range('B11').Value = 'WrappedEther'
$ws.Range('C11').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D11').Value = '1.774.02'
$ws.Range('E11').Value = '  +1.45%  '
$ws.Range('B12').Value = 'TRON'
$ws.Range('C12').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D12').Value = '0.07037'
$ws.Range('E12').Value = '  +1.24%  '
$ws.Range('B13').Value = 'Solana'
$ws.Range('C13').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D13').Value = '15.82'
$ws.Range('E13').Value = '  +4.57%  '
$ws.Range('B14').Value = 'Polygon'
$ws.Range('C14').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D14').Value = '0.6548'
$ws.Range('E14').Value = '  +10.66%  '
$ws.Range('B15').Value = 'Polkadot'
$ws.Range('C15').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D15').Value = '4.497'
$ws.Range('E15').Value = '  +0.43%  '
$ws.Range('B16').Value = 'Litecoin'
$ws.Range('C16').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D16').Value = '78.39'
$ws.Range('E16').Value = '  +2.63%  '
$ws.Range('B17').Value = 'BinanceUSD'
$ws.Range('C17').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D17').Value = '0.9985'
$ws.Range('E17').Value = '  -0.36%  '
$ws.Range('B18').Value = 'Dai'
$ws.Range('C18').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D18').Value = '1.000'
$ws.Range('E18').Value = '  -0.23%  '
$ws.Range('B19').Value = 'WrappedBTC'
$ws.Range('C19').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D19').Value = '26.053.88'
$ws.Range('E19').Value = '  +1.46%  '
$ws.Range('B20').Value = 'Avalanche'
$ws.Range('C20').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D20').Value = '11.69'
$ws.Range('E20').Value = '  +1.12%  '
$ws.Range('B21').Value = 'ShibaInu'
$ws.Range('C21').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D21').Value = '0.000006753'
$ws.Range('E21').Value = '  -0.15%  '
$ws.Range('B22').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C22').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D22').Value = '2.001.12'
$ws.Range('E22').Value = '  +1.60%  '
$ws.Range('B23').Value = 'Uniswap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D23').Value = '4.096'
$ws.Range('E23').Value = '  +1.81%  '
$ws.Range('B24').Value = 'Cosmos'
$ws.Range('C24').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D24').Value = '8.410'
$ws.Range('E24').Value = '  +3.50%  '
$ws.Range('B25').Value = 'Chainlink'
$ws.Range('C25').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D25').Value = '5.195'
$ws.Range('E25').Value = '  +2.27%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').Value = '137.83'
$ws.Range('E26').Value = '  +0.40%  '
$ws.Range('B27').Value = 'Toncoin'
$ws.Range('C27').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D27').Value = '1.483'
$ws.Range('E27').Value = '  -3.00%  '
$ws.Range('B28').Value = 'LidoDAOToken'
$ws.Range('C28').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D28').Value = '1.846'
$ws.Range('E28').Value = '  +1.48%  '
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').Value = '15.21'
$ws.Range('E29').Value = '  +1.79%  '
$ws.Range('B30').Value = 'BitcoinCash'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D30').Value = '102.99'
$ws.Range('E30').Value = '  -0.43%  '
$ws.Range('B31').Value = 'Stellar'
$ws.Range('C31').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D31').Value = '0.08447'
$ws.Range('E31').Value = '  +4.52%  '
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').Value = '3.711'
$ws.Range('E32').Value = '  -1.19%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').Value = '3.441'
$ws.Range('E33').Value = '  -0.29%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').Value = '0.04422'
$ws.Range('E34').Value = '  -1.40%  '
$ws.Range('B35').Value = 'HuobiToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D35').Value = '2.659'
$ws.Range('E35').Value = '  -0.33%  '
$ws.Range('B36').Value = 'ARBITRUM'
$ws.Range('C36').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D36').Value = '1.001'
$ws.Range('E36').Value = '  +2.64%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').Value = '0.6111'
$ws.Range('E37').Value = '  +0.47%  '
$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D38').Value = '2.747'
$ws.Range('E38').Value = '  +3.80%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').Value = '0.01578'
$ws.Range('E39').Value = '  +2.27%  '
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D40').Value = '1.969'
$ws.Range('E40').Value = '  +2.70%  '
$ws.Range('B41').Value = 'PaxDollar'
$ws.Range('C41').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D41').Value = '1.000'
$ws.Range('E41').Value = '  -0.07%  '
$ws.Range('B42').Value = 'Quant'
$ws.Range('C42').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D42').Value = '103.12'
$ws.Range('E42').Value = '  -0.47%  '
$ws.Range('B43').Value = 'TheSandbox'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D43').Value = '0.3896'
$ws.Range('E43').Value = '  +3.09%  '
$ws.Range('B44').Value = 'TrustWalletToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D44').Value = '0.7542'
$ws.Range('E44').Value = '  +3.95%  '
$ws.Range('B45').Value = 'FraxShare'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D45').Value = '4.953'
$ws.Range('E45').Value = '  -2.77%  '
$ws.Range('B46').Value = 'Cronos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D46').Value = '0.05497'
$ws.Range('E46').Value = '  +3.11%  '
$ws.Range('B47').Value = 'Algorand'
$ws.Range('C47').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D47').Value = '0.1122'
$ws.Range('E47').Value = '  +1.47%  '
$ws.Range('B48').Value = 'Aptos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D48').Value = '6.189'
$ws.Range('E48').Value = '  +5.66%  '
$ws.Range('B49').Value = 'Elrond'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D49').Value = '30.20'
$ws.Range('E49').Value = '  +0.64%  '
$ws.Range('B50').Value = 'Aave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D50').Value = '52.93'
$ws.Range('E50').Value = '  +1.24%  '
$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D51').Value = '1.003'
$ws.Range('E51').Value = '  +0.23%  '

Write-Output "Updated cryptos list rows 2-51"